# Table 3 (Normalization of Abundance indicators) formatting fixes:
#  1. Header row: italic column-title runs shrink from 11pt (sz 22) to 10pt (sz 20).
#  2. Data rows: the STOCK-name cell (column 1) switches paragraph alignment
#     from left to right.
#  3. A handful of numeric cells get their text cleaned up (drop stray
#     leading space, pad to 3 decimal places).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Helper: shrink a cell's Range so it no longer includes the trailing
# paragraph-mark / cell-mark characters, trying a few candidate end
# offsets (this runtime's position accounting for the mark pair isn't
# always a flat "-1", so we verify against the expected text and pick
# whichever offset actually lands on it).
# ---------------------------------------------------------------------
function Get-CellContentRange($table, $row, $col, $expectedText) {
    $cell = $table.Cell($row, $col)
    $full = $cell.Range
    $fullText = $full.Text
    $trimmed = $fullText.TrimEnd([char]7, [char]13)
    if ($trimmed -ne $expectedText) {
        throw "Cell ($row,$col): expected [$expectedText] but found [$trimmed]"
    }
    for ($k = 1; $k -le 3; $k++) {
        $candidate = $d.Range($full.Start, $full.End - $k)
        if ($candidate.Text -eq $expectedText) {
            return $candidate
        }
    }
    throw "Cell ($row,$col): could not isolate content range for [$expectedText]"
}

function Set-CellText($table, $row, $col, $oldText, $newText) {
    $rng = Get-CellContentRange $table $row $col $oldText
    $rng.Text = $newText
}

function Set-CellFontSize($table, $row, $col, $size) {
    $cell = $table.Cell($row, $col)
    $full = $cell.Range
    $rng = $d.Range($full.Start, $full.End - 1)
    $rng.Font.Size = $size
}

# ---------------------------------------------------------------------
# 1) Header row (row 1): STOCK, SSBrecent, SSBhistoric, Ftrend, Rtrend,
#    SSBrecent_norm, SSBhistoric_norm, Ftrend_norm, Rtrend_norm,
#    ABUNDANCE -> 10pt
# ---------------------------------------------------------------------
for ($c = 1; $c -le 10; $c++) {
    Set-CellFontSize $t 1 $c 10
}

# ---------------------------------------------------------------------
# 2) STOCK-name column (column 1), data rows 2-12: left -> right align
# ---------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
    $cell = $t.Cell($r, 1)
    $cell.Range.ParagraphFormat.Alignment = 2
}

# ---------------------------------------------------------------------
# 3) Numeric text clean-up
# ---------------------------------------------------------------------
Set-CellText $t 2 4  " 0.007" "0.007"
Set-CellText $t 3 4  " 0.016" "0.016"
Set-CellText $t 3 6  "0"      "0.000"
Set-CellText $t 3 7  "0"      "0.000"
Set-CellText $t 3 10 "0"      "0.000"
Set-CellText $t 4 9  "0"      "0.000"
Set-CellText $t 5 4  " 0.025" "0.025"
Set-CellText $t 5 6  "1"      "1.000"
Set-CellText $t 5 9  "1"      "1.000"
Set-CellText $t 6 4  " 0.000" "0.000"
Set-CellText $t 6 5  "-0.01"  "-0.010"
Set-CellText $t 7 4  " 0.016" "0.016"
Set-CellText $t 8 7  "1"      "1.000"
Set-CellText $t 8 10 "1"      "1.000"
Set-CellText $t 9 4  " 0.043" "0.043"
Set-CellText $t 10 4 " 0.432" "0.432"
Set-CellText $t 11 4 " 0.003" "0.003"
Set-CellText $t 12 4 " 0.011" "0.011"
Set-CellText $t 12 5 "-0.02"  "-0.020"

Write-Host "Done."
